$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

$ws.Range("B1").Value = "Diese Datei wurde erzeugt durch LF-ET 2.2.1 (230325a) und Kommandozeile:"
$ws.Range("B2").Value = '-GenTest "D:/LF/Projekte/rulebased.group/lfet-examples-scope-de/Scopes_01.lfet" -Group "Einfach" -NonExecutableRules "50" -RecommendedTestCases -Protocol -ExportExcel "Scopes_01.testcases.xlsx" -OpenFileOut -OpenDirectoryOut'
$ws.Range("B5").Value = 'Aktuelles Verzeichnis (user.dir): "C:\Program Files\JetBrains\IntelliJ IDEA Community Edition 2022.1.2\jbr\bin"'
$ws.Range("B6").Value = "Benötigte Zeit: 00:00:00.548 (25.03.2023 20:05:14.252 - 25.03.2023 20:05:14.800)"
$ws.Range("B8").Value = 'Entscheidungstabelle: D:\LF\Projekte\rulebased.group\lfet-examples-scope-de\Scopes_01.lfet'
